$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data (genre name, latitude mean, longitude mean) for rows 2-9
$data = @(
    @("historical", 31.96169734751772, 36.09322040969268),
    @("prophecy",   32.04909800329218, 35.39011386927299),
    @("letter",     34.53737994731184, 29.36316048422938),
    @("gospel",     33.77564828008657, 31.73797520014429),
    @("law",        31.65596558106063, 35.77621887525253),
    @("apocalyptic",32.24773232303031, 35.54896866929293),
    @("lyric",      31.05265941960786, 35.46695629150328),
    @("wisdom",     27.77928134285714, 35.48837334285714)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}
